$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "noobject" in G1, matching the style used by the rest of row 1 (F1)
$ws.Range("G1").Value = "noobject"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New column G = PitchThreshold (E) - DiffThreshold (F), for every data row
$ws.Range("G2").Formula = "=E2-F2"
$ws.Range("G3:G13").Formula = "=E3-F3"

# Update the active selection to the newly added column
$ws.Range("G2:G13").Select()

# Drop the stray leftover decorative formatting in the unused H:Q area
# (empty, valueless cells Excel swept away on its own re-save)
$ws.Range("O14:Q14").Clear()
$ws.Range("Q15").Clear()
$ws.Range("H16:Q16").Clear()
$ws.Range("H17:I17").Clear()
$ws.Range("P17:Q17").Clear()
$ws.Range("H18:I18").Clear()
$ws.Range("P18:Q18").Clear()
$ws.Range("H19:Q19").Clear()
$ws.Range("H20:Q20").Clear()
$ws.Range("H21").Clear()
$ws.Range("H22").Clear()
